$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "strapeado" -> "scrapeado" in the Fecha description (B4)
$ws.Range("B4").Value = "una lista de fecha del momento que fue scrapeado en caso de que esos datos mas adelante quisieran utilizar para hacer una comparación"

# Row 4: increase height, underline the label cell (A4)
$ws.Range("A4").Font.Underline = $true
$ws.Rows.Item(4).RowHeight = 35.4

# New blank, underlined placeholder cells added further down the sheet
$ws.Range("B9").Font.Underline = $true
$ws.Range("B9").Font.Size = 10
$ws.Rows.Item(9).RowHeight = 15.75

$ws.Range("B10").Font.Underline = $true
$ws.Range("B10").Font.Size = 10
$ws.Rows.Item(10).RowHeight = 15.75

$ws.Range("B13").Font.Underline = $true
$ws.Range("B13").Font.Size = 10
$ws.Rows.Item(13).RowHeight = 15.75

# Update selection to match the author's last active cell
$null = $ws.Range("B8").Select()
